# Daily attendance processing - 2026-01-17 09:03:48
# Replace "dnasr281@gmail.com, System" with "System, dnasr281@gmail.com"
# in column G (Recorded By) wherever that exact text occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$colG = $ws.Range("G1:G$lastRow")

# xlWhole = 1 (match entire cell contents, not a substring)
$colG.Replace(
    "dnasr281@gmail.com, System",
    "System, dnasr281@gmail.com",
    1,
    1,
    $false,
    $false,
    $false
)
